# Update low input values for other language (per-row run statistics,
# columns run_time, num_deaths, max_er, and iter 0..19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value2 = 1.123036861419678
$ws.Cells.Item(2, 4).Value2 = 0
$ws.Cells.Item(2, 5).Value2 = 2219.252388005494
$ws.Cells.Item(2, 6).Value2 = 0.1389755513198633
$ws.Cells.Item(2, 7).Value2 = 0.09847495715240451
$ws.Cells.Item(2, 8).Value2 = 0.07598346652009627
$ws.Cells.Item(2, 9).Value2 = 0.06932536610657877
$ws.Cells.Item(2, 10).Value2 = 0.06545977253632554
$ws.Cells.Item(2, 11).Value2 = 0.05837246123804143
$ws.Cells.Item(2, 12).Value2 = 0.05458221121676798
$ws.Cells.Item(2, 13).Value2 = 0.05011680591399206
$ws.Cells.Item(2, 14).Value2 = 0.04906164072771402
$ws.Cells.Item(2, 15).Value2 = 0.04725081958330594
$ws.Cells.Item(2, 16).Value2 = 0.04684432312843756
$ws.Cells.Item(2, 17).Value2 = 0.04574533391165668
$ws.Cells.Item(2, 18).Value2 = 0.04574533391165668
$ws.Cells.Item(2, 19).Value2 = 0.04514911740432084
$ws.Cells.Item(2, 20).Value2 = 0.04426259194101713
$ws.Cells.Item(2, 21).Value2 = 0.04426259194101713
$ws.Cells.Item(2, 22).Value2 = 0.04378265349141725
$ws.Cells.Item(2, 23).Value2 = 0.04350713578155373
$ws.Cells.Item(2, 24).Value2 = 0.04346465259067465
$ws.Cells.Item(2, 25).Value2 = 0.04326028046794334

$ws.Cells.Item(3, 3).Value2 = 1.199986219406128
$ws.Cells.Item(3, 4).Value2 = 0
$ws.Cells.Item(3, 5).Value2 = 2337.368060991624
$ws.Cells.Item(3, 6).Value2 = 0.1257266616179585
$ws.Cells.Item(3, 7).Value2 = 0.1038741844516984
$ws.Cells.Item(3, 8).Value2 = 0.09056011235244439
$ws.Cells.Item(3, 9).Value2 = 0.0801525765668521
$ws.Cells.Item(3, 10).Value2 = 0.06851028614192411
$ws.Cells.Item(3, 11).Value2 = 0.06480875698685114
$ws.Cells.Item(3, 12).Value2 = 0.05740606483237725
$ws.Cells.Item(3, 13).Value2 = 0.05448862348109607
$ws.Cells.Item(3, 14).Value2 = 0.05205527252207066
$ws.Cells.Item(3, 15).Value2 = 0.05184030110026658
$ws.Cells.Item(3, 16).Value2 = 0.0499548244517927
$ws.Cells.Item(3, 17).Value2 = 0.04868184712049774
$ws.Cells.Item(3, 18).Value2 = 0.04827873861326929
$ws.Cells.Item(3, 19).Value2 = 0.04706791988510936
$ws.Cells.Item(3, 20).Value2 = 0.04682467396382849
$ws.Cells.Item(3, 21).Value2 = 0.04625693658023868
$ws.Cells.Item(3, 22).Value2 = 0.04615508989898263
$ws.Cells.Item(3, 23).Value2 = 0.04581041677222375
$ws.Cells.Item(3, 24).Value2 = 0.04581041677222375
$ws.Cells.Item(3, 25).Value2 = 0.04556273023375484

$ws.Cells.Item(4, 3).Value2 = 1.375
$ws.Cells.Item(4, 4).Value2 = 0
$ws.Cells.Item(4, 5).Value2 = 2252.771068130824
$ws.Cells.Item(4, 6).Value2 = 0.1512467526428288
$ws.Cells.Item(4, 7).Value2 = 0.09499009377469206
$ws.Cells.Item(4, 8).Value2 = 0.08043805932696277
$ws.Cells.Item(4, 9).Value2 = 0.06926126580168042
$ws.Cells.Item(4, 10).Value2 = 0.06453962817486562
$ws.Cells.Item(4, 11).Value2 = 0.06134019609932842
$ws.Cells.Item(4, 12).Value2 = 0.05600116185411672
$ws.Cells.Item(4, 13).Value2 = 0.0536923913820142
$ws.Cells.Item(4, 14).Value2 = 0.05129270945631941
$ws.Cells.Item(4, 15).Value2 = 0.04913494404738877
$ws.Cells.Item(4, 16).Value2 = 0.04718225232736681
$ws.Cells.Item(4, 17).Value2 = 0.04718225232736681
$ws.Cells.Item(4, 18).Value2 = 0.04646838305900025
$ws.Cells.Item(4, 19).Value2 = 0.04582499379796614
$ws.Cells.Item(4, 20).Value2 = 0.04500652688431959
$ws.Cells.Item(4, 21).Value2 = 0.04422385703910584
$ws.Cells.Item(4, 22).Value2 = 0.04422385703910584
$ws.Cells.Item(4, 23).Value2 = 0.04419953905642796
$ws.Cells.Item(4, 24).Value2 = 0.04403824477800929
$ws.Cells.Item(4, 25).Value2 = 0.04391366604543515

$ws.Cells.Item(5, 3).Value2 = 1.201002359390259
$ws.Cells.Item(5, 4).Value2 = 0
$ws.Cells.Item(5, 5).Value2 = 2333.545590465834
$ws.Cells.Item(5, 6).Value2 = 0.1423934733437842
$ws.Cells.Item(5, 7).Value2 = 0.1075409253115023
$ws.Cells.Item(5, 8).Value2 = 0.07538189070786012
$ws.Cells.Item(5, 9).Value2 = 0.07323341581521252
$ws.Cells.Item(5, 10).Value2 = 0.06833327064604522
$ws.Cells.Item(5, 11).Value2 = 0.06197989599098019
$ws.Cells.Item(5, 12).Value2 = 0.05878651184678445
$ws.Cells.Item(5, 13).Value2 = 0.05671742193459618
$ws.Cells.Item(5, 14).Value2 = 0.05465262634704669
$ws.Cells.Item(5, 15).Value2 = 0.05200187263809779
$ws.Cells.Item(5, 16).Value2 = 0.0505353039414032
$ws.Cells.Item(5, 17).Value2 = 0.04971022953939917
$ws.Cells.Item(5, 18).Value2 = 0.04865619377463083
$ws.Cells.Item(5, 19).Value2 = 0.04778672807836091
$ws.Cells.Item(5, 20).Value2 = 0.04704423263782968
$ws.Cells.Item(5, 21).Value2 = 0.04657777666703193
$ws.Cells.Item(5, 22).Value2 = 0.0464464998713536
$ws.Cells.Item(5, 23).Value2 = 0.04597240442965984
$ws.Cells.Item(5, 24).Value2 = 0.04569273843267307
$ws.Cells.Item(5, 25).Value2 = 0.04548821813773554

$ws.Cells.Item(6, 3).Value2 = 1.210992336273193
$ws.Cells.Item(6, 4).Value2 = 0
$ws.Cells.Item(6, 5).Value2 = 2270.640212300119
$ws.Cells.Item(6, 6).Value2 = 0.1443284787734793
$ws.Cells.Item(6, 7).Value2 = 0.1032272779495481
$ws.Cells.Item(6, 8).Value2 = 0.08435666994202849
$ws.Cells.Item(6, 9).Value2 = 0.06500711861855689
$ws.Cells.Item(6, 10).Value2 = 0.0646258356207151
$ws.Cells.Item(6, 11).Value2 = 0.06056557031379652
$ws.Cells.Item(6, 12).Value2 = 0.05658659393690224
$ws.Cells.Item(6, 13).Value2 = 0.05383554837980897
$ws.Cells.Item(6, 14).Value2 = 0.05146342812749059
$ws.Cells.Item(6, 15).Value2 = 0.04965256077965503
$ws.Cells.Item(6, 16).Value2 = 0.04830972570964767
$ws.Cells.Item(6, 17).Value2 = 0.04771703752352512
$ws.Cells.Item(6, 18).Value2 = 0.04618812948942698
$ws.Cells.Item(6, 19).Value2 = 0.04549729879815737
$ws.Cells.Item(6, 20).Value2 = 0.04486967319451038
$ws.Cells.Item(6, 21).Value2 = 0.04456589430511874
$ws.Cells.Item(6, 22).Value2 = 0.04456589430511874
$ws.Cells.Item(6, 23).Value2 = 0.04441347299832551
$ws.Cells.Item(6, 24).Value2 = 0.04440283981202425
$ws.Cells.Item(6, 25).Value2 = 0.04426199244249743

$ws.Cells.Item(7, 3).Value2 = 1.224040985107422
$ws.Cells.Item(7, 4).Value2 = 0
$ws.Cells.Item(7, 5).Value2 = 2189.848677939339
$ws.Cells.Item(7, 6).Value2 = 0.1197607356035581
$ws.Cells.Item(7, 7).Value2 = 0.1015236962583557
$ws.Cells.Item(7, 8).Value2 = 0.07777786296457963
$ws.Cells.Item(7, 9).Value2 = 0.06614975410877255
$ws.Cells.Item(7, 10).Value2 = 0.06606662102624862
$ws.Cells.Item(7, 11).Value2 = 0.05806128172266541
$ws.Cells.Item(7, 12).Value2 = 0.05327493690625272
$ws.Cells.Item(7, 13).Value2 = 0.05173417330031371
$ws.Cells.Item(7, 14).Value2 = 0.05023942492992341
$ws.Cells.Item(7, 15).Value2 = 0.04841136220950151
$ws.Cells.Item(7, 16).Value2 = 0.0464429508175998
$ws.Cells.Item(7, 17).Value2 = 0.04602550320020914
$ws.Cells.Item(7, 18).Value2 = 0.04451892323800887
$ws.Cells.Item(7, 19).Value2 = 0.04357704162935757
$ws.Cells.Item(7, 20).Value2 = 0.04339369745673335
$ws.Cells.Item(7, 21).Value2 = 0.04322448460494735
$ws.Cells.Item(7, 22).Value2 = 0.0431534866039608
$ws.Cells.Item(7, 23).Value2 = 0.04288978337529514
$ws.Cells.Item(7, 24).Value2 = 0.04271844458656645
$ws.Cells.Item(7, 25).Value2 = 0.04268710873176098

$ws.Cells.Item(8, 3).Value2 = 1.131015539169312
$ws.Cells.Item(8, 4).Value2 = 0
$ws.Cells.Item(8, 5).Value2 = 2250.797205590547
$ws.Cells.Item(8, 6).Value2 = 0.1286959795927647
$ws.Cells.Item(8, 7).Value2 = 0.0977665817028233
$ws.Cells.Item(8, 8).Value2 = 0.08059309783373239
$ws.Cells.Item(8, 9).Value2 = 0.07101441649348225
$ws.Cells.Item(8, 10).Value2 = 0.06372399223194741
$ws.Cells.Item(8, 11).Value2 = 0.05889815534424102
$ws.Cells.Item(8, 12).Value2 = 0.05577744783027858
$ws.Cells.Item(8, 13).Value2 = 0.05352636850654065
$ws.Cells.Item(8, 14).Value2 = 0.0519981608777314
$ws.Cells.Item(8, 15).Value2 = 0.04959292942520523
$ws.Cells.Item(8, 16).Value2 = 0.0488710827049555
$ws.Cells.Item(8, 17).Value2 = 0.0472368796250404
$ws.Cells.Item(8, 18).Value2 = 0.04645319620514955
$ws.Cells.Item(8, 19).Value2 = 0.04583963097529657
$ws.Cells.Item(8, 20).Value2 = 0.0452765267161289
$ws.Cells.Item(8, 21).Value2 = 0.04455797305124105
$ws.Cells.Item(8, 22).Value2 = 0.04455797305124105
$ws.Cells.Item(8, 23).Value2 = 0.0442198084711523
$ws.Cells.Item(8, 24).Value2 = 0.04414803310554626
$ws.Cells.Item(8, 25).Value2 = 0.04387518919279817

$ws.Cells.Item(9, 3).Value2 = 1.225028514862061
$ws.Cells.Item(9, 4).Value2 = 0
$ws.Cells.Item(9, 5).Value2 = 2174.861587870793
$ws.Cells.Item(9, 6).Value2 = 0.1440219836491839
$ws.Cells.Item(9, 7).Value2 = 0.08769392331245404
$ws.Cells.Item(9, 8).Value2 = 0.08173941032139932
$ws.Cells.Item(9, 9).Value2 = 0.06995847025039584
$ws.Cells.Item(9, 10).Value2 = 0.06523278814191104
$ws.Cells.Item(9, 11).Value2 = 0.05724987647205408
$ws.Cells.Item(9, 12).Value2 = 0.05413875727052933
$ws.Cells.Item(9, 13).Value2 = 0.05080358988782813
$ws.Cells.Item(9, 14).Value2 = 0.04857006722022774
$ws.Cells.Item(9, 15).Value2 = 0.04714881360254937
$ws.Cells.Item(9, 16).Value2 = 0.04655628504262817
$ws.Cells.Item(9, 17).Value2 = 0.04504145733464982
$ws.Cells.Item(9, 18).Value2 = 0.04483455068931814
$ws.Cells.Item(9, 19).Value2 = 0.04356978484037112
$ws.Cells.Item(9, 20).Value2 = 0.04337210665841894
$ws.Cells.Item(9, 21).Value2 = 0.04335465133373347
$ws.Cells.Item(9, 22).Value2 = 0.04289151180333093
$ws.Cells.Item(9, 23).Value2 = 0.04283855131309471
$ws.Cells.Item(9, 24).Value2 = 0.04246047862279198
$ws.Cells.Item(9, 25).Value2 = 0.04239496272652618

$ws.Cells.Item(10, 3).Value2 = 1.148012399673462
$ws.Cells.Item(10, 4).Value2 = 0
$ws.Cells.Item(10, 5).Value2 = 2310.580258079313
$ws.Cells.Item(10, 6).Value2 = 0.1294850815401244
$ws.Cells.Item(10, 7).Value2 = 0.10277947172215
$ws.Cells.Item(10, 8).Value2 = 0.08177881137298992
$ws.Cells.Item(10, 9).Value2 = 0.07476636083785003
$ws.Cells.Item(10, 10).Value2 = 0.0684105838079693
$ws.Cells.Item(10, 11).Value2 = 0.06490523580969185
$ws.Cells.Item(10, 12).Value2 = 0.06020386807853812
$ws.Cells.Item(10, 13).Value2 = 0.05545533534790918
$ws.Cells.Item(10, 14).Value2 = 0.05275036894998396
$ws.Cells.Item(10, 15).Value2 = 0.05074568728088612
$ws.Cells.Item(10, 16).Value2 = 0.04990185711241785
$ws.Cells.Item(10, 17).Value2 = 0.04841707562471678
$ws.Cells.Item(10, 18).Value2 = 0.04760329344006436
$ws.Cells.Item(10, 19).Value2 = 0.04645888451947912
$ws.Cells.Item(10, 20).Value2 = 0.04637280547695979
$ws.Cells.Item(10, 21).Value2 = 0.04578080689349801
$ws.Cells.Item(10, 22).Value2 = 0.04576673919761902
$ws.Cells.Item(10, 23).Value2 = 0.04543614859818033
$ws.Cells.Item(10, 24).Value2 = 0.04504055083975268
$ws.Cells.Item(10, 25).Value2 = 0.04504055083975268

$ws.Cells.Item(11, 3).Value2 = 1.304998874664307
$ws.Cells.Item(11, 4).Value2 = 0
$ws.Cells.Item(11, 5).Value2 = 2311.784103735229
$ws.Cells.Item(11, 6).Value2 = 0.1336344263736033
$ws.Cells.Item(11, 7).Value2 = 0.1011439233109617
$ws.Cells.Item(11, 8).Value2 = 0.08304347712868669
$ws.Cells.Item(11, 9).Value2 = 0.07073236975819963
$ws.Cells.Item(11, 10).Value2 = 0.06343784245802858
$ws.Cells.Item(11, 11).Value2 = 0.06125771008239581
$ws.Cells.Item(11, 12).Value2 = 0.05705560883509756
$ws.Cells.Item(11, 13).Value2 = 0.05395228005449909
$ws.Cells.Item(11, 14).Value2 = 0.05233196665173879
$ws.Cells.Item(11, 15).Value2 = 0.04934077563711561
$ws.Cells.Item(11, 16).Value2 = 0.04934077563711561
$ws.Cells.Item(11, 17).Value2 = 0.0477060610864339
$ws.Cells.Item(11, 18).Value2 = 0.04753298000287777
$ws.Cells.Item(11, 19).Value2 = 0.04725379448153589
$ws.Cells.Item(11, 20).Value2 = 0.04712633929478422
$ws.Cells.Item(11, 21).Value2 = 0.04617653456264369
$ws.Cells.Item(11, 22).Value2 = 0.04572082181859721
$ws.Cells.Item(11, 23).Value2 = 0.04546716810525685
$ws.Cells.Item(11, 24).Value2 = 0.04506401761667112
$ws.Cells.Item(11, 25).Value2 = 0.04506401761667112

